# Apply updated cryptocurrency market data to Sheet1.
# (commit: "Updated cryptos list ... with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.579.52"
$ws.Range("D3").Value = "2.285.75"
$ws.Range("E4").Value = "  -0.03%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "304.71"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.59%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "95.52"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -2.38%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.503"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -3.39%  "
$ws.Range("E8").Value = "  +0.02%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.496"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -3.46%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "34.95"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -3.14%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0782"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.86%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "18.20"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +2.98%  "
$ws.Range("E13").Value = "  +1.05%  "
$ws.Range("E14").Value = "  -2.54%  "
$ws.Range("D15").Value = "2.644.80"
$ws.Range("E15").Value = "  -0.46%  "
$ws.Range("D16").Value = "2.285.28"
$ws.Range("E16").Value = "  -1.02%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.773"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -1.59%  "
$ws.Range("D18").Value = "42.494.64"
$ws.Range("E18").Value = "  -1.07%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "12.77"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "0.0₃0890"
$ws.Range("E20").Value = "  -2.70%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "5.98"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -2.20%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "66.75"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -3.39%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "235.21"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.90%  "
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("E27").Value = "  +0.06%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "24.88"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.44%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "165.80"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  +0.34%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "8.97"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -1.42%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "32.54"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.73%  "
$ws.Range("E33").Value = "  +0.00%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.67"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -1.36%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "4.92"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -3.04%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "17.49"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -2.11%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.39"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.81%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.0685"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.50%  "
$ws.Range("E39").Value = "  -1.11%  "
$ws.Range("E40").Value = "  -2.31%  "
$ws.Range("E41").Value = "  -1.88%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "2.67"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -3.34%  "
$ws.Range("D43").Value = "1.990.56"
$ws.Range("E43").Value = "  -0.86%  "
$ws.Range("E44").Value = "  -2.94%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "9.97"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -2.86%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "17.92"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +3.12%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "2.03"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -9.27%  "
$ws.Range("E48").Value = "  -2.34%  "
$ws.Range("E49").Value = "  +8.49%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "53.28"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.46%  "
$ws.Range("D51").Value = "2.511.86"
$ws.Range("E51").Value = "  -0.40%  "
